# Apply the crypto-price refresh captured in the commit diff.
# D-column cells store plain-text numbers (e.g. "1.000", "30.130.99") in the
# original workbook (t="inlineStr"). Excel's Range.Value setter auto-coerces
# numeric-looking strings to real numbers, which would silently drop the
# formatting (e.g. "1.000" -> 1). Prefixing with a leading apostrophe forces
# Excel to keep/store the value as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''30.130.99'
$ws.Range("E2").Value = '  +5.51%  '

# Row 3
$ws.Range("D3").Value = '''1.922.45'
$ws.Range("E3").Value = '  +2.37%  '

# Row 4
$ws.Range("E4").Value = '  -1.08%  '

# Row 5
$ws.Range("D5").Value = '''327.74'
$ws.Range("E5").Value = '  +3.65%  '

# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  -0.87%  '

# Row 7
$ws.Range("D7").Value = '''0.5164'
$ws.Range("E7").Value = '  +1.36%  '

# Row 8
$ws.Range("D8").Value = '''0.4005'
$ws.Range("E8").Value = '  +2.56%  '

# Row 9
$ws.Range("D9").Value = '''0.08460'
$ws.Range("E9").Value = '  +0.72%  '

# Row 10
$ws.Range("D10").Value = '''42.88'
$ws.Range("E10").Value = '  +2.34%  '

# Row 11
$ws.Range("D11").Value = '''1.121'
$ws.Range("E11").Value = '  +1.58%  '

# Row 12
$ws.Range("D12").Value = '''21.28'
$ws.Range("E12").Value = '  +4.21%  '

# Row 13
$ws.Range("D13").Value = '''6.339'
$ws.Range("E13").Value = '  +1.78%  '

# Row 14
$ws.Range("D14").Value = '''1.917.67'
$ws.Range("E14").Value = '  +2.22%  '

# Row 15
$ws.Range("D15").Value = '''7.340'
$ws.Range("E15").Value = '  +1.27%  '

# Row 16
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  -1.11%  '

# Row 17
$ws.Range("D17").Value = '''95.04'
$ws.Range("E17").Value = '  +4.05%  '

# Row 18
$ws.Range("D18").Value = '''0.00001116'
$ws.Range("E18").Value = '  +1.06%  '

# Row 19
$ws.Range("D19").Value = '''0.06730'
$ws.Range("E19").Value = '  +0.16%  '

# Row 20
$ws.Range("D20").Value = '''18.06'
$ws.Range("E20").Value = '  +1.94%  '

# Row 21
$ws.Range("E21").Value = '  -0.81%  '

# Row 22
$ws.Range("D22").Value = '''6.068'
$ws.Range("E22").Value = '  +2.27%  '

# Row 23
$ws.Range("D23").Value = '''30.132.31'

# Row 24
$ws.Range("E24").Value = '  +1.02%  '

# Row 25
$ws.Range("D25").Value = '''2.205'
$ws.Range("E25").Value = '  -1.44%  '

# Row 26
$ws.Range("D26").Value = '''2.139.43'
$ws.Range("E26").Value = '  +2.49%  '

# Row 27
$ws.Range("D27").Value = '''160.95'
$ws.Range("E27").Value = '  -0.43%  '

# Row 28
$ws.Range("D28").Value = '''20.97'
$ws.Range("E28").Value = '  +1.42%  '

# Row 29
$ws.Range("D29").Value = '''2.462'
$ws.Range("E29").Value = '  +4.71%  '

# Row 30
$ws.Range("D30").Value = '''128.82'
$ws.Range("E30").Value = '  +2.19%  '

# Row 31
$ws.Range("E31").Value = '  +3.25%  '

# Row 32
$ws.Range("D32").Value = '''0.1061'
$ws.Range("E32").Value = '  +1.43%  '

# Row 33
$ws.Range("D33").Value = '''6.076'
$ws.Range("E33").Value = '  +4.92%  '

# Row 34
$ws.Range("D34").Value = '''3.659'
$ws.Range("E34").Value = '  +1.33%  '

# Row 35
$ws.Range("D35").Value = '''0.02509'
$ws.Range("E35").Value = '  +1.93%  '

# Row 36
$ws.Range("D36").Value = '''0.06595'
$ws.Range("E36").Value = '  +0.69%  '

# Row 37
$ws.Range("D37").Value = '''0.2227'
$ws.Range("E37").Value = '  +2.93%  '

# Row 38
$ws.Range("D38").Value = '''1.238'
$ws.Range("E38").Value = '  +3.56%  '

# Row 39
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").Value = '''5.212'
$ws.Range("E39").Value = '  +2.82%  '

# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''8.998'
$ws.Range("E40").Value = '  +1.66%  '

# Row 41
$ws.Range("E41").Value = '  +1.91%  '

# Row 42
$ws.Range("D42").Value = '''1.241'
$ws.Range("E42").Value = '  -0.68%  '

# Row 43
$ws.Range("E43").Value = '  +2.36%  '

# Row 44
$ws.Range("D44").Value = '''0.6138'
$ws.Range("E44").Value = '  +1.62%  '

# Row 45
$ws.Range("D45").Value = '''13.17'
$ws.Range("E45").Value = '  +1.16%  '

# Row 46
$ws.Range("D46").Value = '''3.764'
$ws.Range("E46").Value = '  +1.93%  '

# Row 47
$ws.Range("D47").Value = '''2.054'
$ws.Range("E47").Value = '  +2.16%  '

# Row 48
$ws.Range("D48").Value = '''125.96'
$ws.Range("E48").Value = '  +3.14%  '

# Row 49
$ws.Range("D49").Value = '''1.242'
$ws.Range("E49").Value = '  +1.94%  '

# Row 50
$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = '''1.156'
$ws.Range("E50").Value = '  +0.68%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''79.33'
$ws.Range("E51").Value = '  +3.18%  '
